# Memory Benchmark = 255,320,430 bytes
#
# 1) Rename "Sheet1" -> "VS Heap Tool"
# 2) Add a new sheet "Heap Report from Test" right after it, and populate it
#    with the heap-memory benchmark report (Date/Platform/Build/Branch/
#    Duration/Heap Memory Test 1-3/Description/Observations).
# 3) Clean up the view state on "VS Heap Tool" (selection -> A1:E1) and make
#    "Heap Report from Test" the active tab at 150% zoom.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "VS Heap Tool"

# Insert the new sheet immediately after "VS Heap Tool"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Heap Report from Test"

# --- Header row (row 1) -----------------------------------------------
# Write cells in the same order the original workbook's shared-string
# table picked them up in, so new unique strings land at matching indices.
$ws2.Range("A1").Value = "Date"

$ws2.Range("D1").NumberFormat = "#,##0"
$ws2.Range("D1").Value = "Branch"

$ws2.Range("E1").Value = "Duration"

$ws2.Range("F1").NumberFormat = "#,##0"
$ws2.Range("F1").Value = "Heap Memory Test 1"

$ws2.Range("G1").NumberFormat = "#,##0"
$ws2.Range("G1").Value = "Heap Memory Test 2"

$ws2.Range("H1").NumberFormat = "#,##0"
$ws2.Range("H1").Value = "Heap Memory Test 3"

$ws2.Range("I1").WrapText = $true
$ws2.Range("I1").Value = "Description"

$ws2.Range("J1").WrapText = $true
$ws2.Range("J1").Value = "Observations"

# --- Data row (row 2) ---------------------------------------------------
# Reuse the date/time number format already in the workbook (same format
# as column A on "VS Heap Tool") by copying its formatting over.
$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Value = 43409.84375

$ws2.Range("D2").Value = "Research_Memory_Management_Options"

$ws2.Range("B1").Value = "Platform"
$ws2.Range("B2").Value = "PC"

$ws2.Range("C1").Value = "Build"
$ws2.Range("C2").Value = "Release"

$ws2.Range("F2").NumberFormat = "#,##0"
$ws2.Range("F2").Value = 255320430

$ws2.Range("G2").NumberFormat = "#,##0"
$ws2.Range("G2").Value = 255320430

$ws2.Range("H2").NumberFormat = "#,##0"
$ws2.Range("H2").Value = 255320430

# --- Column widths (best-fit-ish) --------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 14
$ws2.Columns.Item(4).ColumnWidth = 35.333333333333336
$ws2.Columns.Item(6).ColumnWidth = 17.333333333333332
$ws2.Columns.Item(7).ColumnWidth = 17.333333333333332
$ws2.Columns.Item(8).ColumnWidth = 17.333333333333332
$ws2.Columns.Item(9).ColumnWidth = 49.833333333333336
$ws2.Columns.Item(10).ColumnWidth = 49.833333333333336

$ws2.PageSetup.Orientation = 1

# --- View state ----------------------------------------------------------
# "VS Heap Tool" keeps a plain A1:E1 selection and is no longer the active
# tab / no longer tab-selected.
$ws1.Range("A1:E1").Select() | Out-Null

# "Heap Report from Test" becomes the active sheet, zoomed to 150%.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 150
